$wb = $excel.ActiveWorkbook

# --- 1. Insert new sheet "2022-Q1" right after "2021-Q4" (before "总计") ---
$afterSheet = $wb.Worksheets.Item("2021-Q4")
$newSheet = $wb.Worksheets.Add($null, $afterSheet)
$newSheet.Name = "2022-Q1"

# Header row
$newSheet.Range("B1").Value = "基金代码"
$newSheet.Range("C1").Value = "基金名称"
$newSheet.Range("D1").Value = "基金规模"
$newSheet.Range("E1").Value = "股票总仓位"
$newSheet.Range("F1").Value = "仓位占比"
$newSheet.Range("G1").Value = "持有市值(亿元)"
$newSheet.Range("H1").Value = "仓位排名"

# Fund code column (B) must stay text so leading zeros survive
$newSheet.Range("B2:B7").NumberFormat = "@"
# Numeric-looking text columns (D,E,F,G) are stored as text in this workbook
$newSheet.Range("D2:D7").NumberFormat = "@"
$newSheet.Range("E2:E7").NumberFormat = "@"
$newSheet.Range("F2:F7").NumberFormat = "@"
$newSheet.Range("G2:G7").NumberFormat = "@"

# Data rows
$newSheet.Range("A2").Value = 0
$newSheet.Range("B2").Value = "012368"
$newSheet.Range("C2").Value = "摩根士丹利华鑫优享臻选六个月持有期混合型证券投资基金A"
$newSheet.Range("D2").Value = "5.76"
$newSheet.Range("E2").Value = "93.78"
$newSheet.Range("F2").Value = "10.06"
$newSheet.Range("G2").Value = "0.5795"
$newSheet.Range("H2").Value = 1

$newSheet.Range("A3").Value = 1
$newSheet.Range("B3").Value = "233006"
$newSheet.Range("C3").Value = "大摩领先优势混合"
$newSheet.Range("D3").Value = "4.12"
$newSheet.Range("E3").Value = "94.42"
$newSheet.Range("F3").Value = "10.15"
$newSheet.Range("G3").Value = "0.4182"
$newSheet.Range("H3").Value = 1

$newSheet.Range("A4").Value = 2
$newSheet.Range("B4").Value = "000309"
$newSheet.Range("C4").Value = "大摩品质生活精选股票"
$newSheet.Range("D4").Value = "4.36"
$newSheet.Range("E4").Value = "94.17"
$newSheet.Range("F4").Value = "8.76"
$newSheet.Range("G4").Value = "0.3819"
$newSheet.Range("H4").Value = 1

$newSheet.Range("A5").Value = 3
$newSheet.Range("B5").Value = "010322"
$newSheet.Range("C5").Value = "摩根士丹利华鑫新兴产业股票"
$newSheet.Range("D5").Value = "2.41"
$newSheet.Range("E5").Value = "94.11"
$newSheet.Range("F5").Value = "10.35"
$newSheet.Range("G5").Value = "0.2494"
$newSheet.Range("H5").Value = 1

$newSheet.Range("A6").Value = 4
$newSheet.Range("B6").Value = "002707"
$newSheet.Range("C6").Value = "摩根士丹利华鑫科技领先灵活配置混合"
$newSheet.Range("D6").Value = "2.27"
$newSheet.Range("E6").Value = "93.05"
$newSheet.Range("F6").Value = "5.03"
$newSheet.Range("G6").Value = "0.1142"
$newSheet.Range("H6").Value = 7

$newSheet.Range("A7").Value = 5
$newSheet.Range("B7").Value = "012369"
$newSheet.Range("C7").Value = "摩根士丹利华鑫优享臻选六个月持有期混合型证券投资基金C"
$newSheet.Range("D7").Value = "0.40"
$newSheet.Range("E7").Value = "93.78"
$newSheet.Range("F7").Value = "10.06"
$newSheet.Range("G7").Value = "0.0402"
$newSheet.Range("H7").Value = 1

# Apply header style (bold/centered/bordered) matching the other sheets.
# Applied per-cell: doing it over a multi-cell range at once does not
# reliably persist the border in this runtime.
foreach ($addr in @("B1","C1","D1","E1","F1","G1","H1","A2","A3","A4","A5","A6","A7")) {
    $cell = $newSheet.Range($addr)
    $cell.Font.Bold = $true
    $cell.HorizontalAlignment = -4108
    $cell.VerticalAlignment = -4160
    $cell.Borders.LineStyle = 1
}

# --- 2. Insert a new top data row in "总计" sheet for 2022-Q1 ---
$totalSheet = $wb.Worksheets.Item("总计")
$totalSheet.Rows.Item(2).Insert()

$totalSheet.Range("A2").Value = 0
$totalSheet.Range("B2").Value = "2022-Q1"
$totalSheet.Range("C2").Value = 6
$totalSheet.Range("D2").Value = 1.78

# Re-number the index column (A) for the rows that shifted down
$totalSheet.Range("A3").Value = 1
$totalSheet.Range("A4").Value = 2
$totalSheet.Range("A5").Value = 3
$totalSheet.Range("A6").Value = 4
$totalSheet.Range("A7").Value = 5
